# This script "repulls" the dSF (column F) data for the montgomery_jordan
# workbook, updating the F-column cell values to reflect the freshly
# pulled/recalculated figures, leaving all other cells untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dsfUpdates = @{
    2  = -1
    3  = 6
    4  = -1
    5  = 2
    6  = 5
    7  = 4
    8  = 1
    9  = 2
    10 = 1
    11 = 7
    14 = 5
    15 = 1
    16 = 3
    17 = -2
    18 = 2
    20 = -6
    21 = -1
    22 = 4
    23 = 2
    24 = 2
    25 = 5
    26 = 3
    27 = 5
    29 = 1
    30 = -1
    31 = -1
    32 = -4
    33 = 1
    34 = -1
    36 = -5
    37 = -1
    38 = 6
    39 = 6
    40 = 4
    42 = 3
    43 = -3
}

foreach ($row in $dsfUpdates.Keys) {
    $ws.Range("F$row").Value = $dsfUpdates[$row]
}
